$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap columns A and B for header + existing data rows (1-26) ---
for ($r = 1; $r -le 26; $r++) {
    $a = $ws.Cells.Item($r, 1).Value()
    $b = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($r, 1).Value = $b
    $ws.Cells.Item($r, 2).Value = $a
}

# --- Append new rows 27-38 ---
$newRows = @(
    @("Jimmy Choo", "Urban Hero Gold ", "woda perfumowana dla mężczyzn", "50 ml ", 133.5, "https://www.notino.pl/jimmy-choo/urban-hero-gold-woda-perfumowana-dla-mezczyzn/", "10.03.2023"),
    @("Jimmy Choo", "Urban Hero Gold ", "woda perfumowana dla mężczyzn", "100 ml ", 217.5, "https://www.notino.pl/jimmy-choo/urban-hero-gold-woda-perfumowana-dla-mezczyzn/", "04.04.2023"),
    @("JOOP!", "Homme ", "woda toaletowa dla mężczyzn", "200 ml ", 196, "https://www.notino.pl/joop/homme-woda-toaletowa-dla-mczyzn/p-402506/?gclid=EAIaIQobChMInZmGhfyU_gIVg9eyCh16ygIcEAQYASABEgJDUPD_BwE", "06.04.2023"),
    @("JOOP!", "Homme ", "woda toaletowa dla mężczyzn", "200 ml ", 196, "https://www.notino.pl/joop/homme-woda-toaletowa-dla-mczyzn/p-402506/?gclid=EAIaIQobChMInZmGhfyU_gIVg9eyCh16ygIcEAQYASABEgJDUPD_BwE", "06.04.2023"),
    @("Armani", "Sì ", "woda perfumowana dla kobiet", "50 ml ", 273.28, "https://www.notino.pl/armani/si-woda-perfumowana-dla-kobiet/", "07.04.2023"),
    @("Armani", "Sì ", "woda perfumowana dla kobiet", "50 ml ", 375, "https://www.notino.pl/armani/si-woda-perfumowana-dla-kobiet/", "10.04.2023"),
    @("Armani", "Sì ", "woda perfumowana dla kobiet", "50 ml ", 375, "https://www.notino.pl/armani/si-woda-perfumowana-dla-kobiet/", "10.04.2023"),
    @("Paco Rabanne", "Invictus Victory Elixir ", "perfumy dla mężczyzn", "200 ml ", 885, "https://www.notino.pl/paco-rabanne/invictus-victory-elixir-perfumy-dla-mezczyzn/p-16168428/", "11.04.2023"),
    @("Paco Rabanne", "Invictus Victory Elixir ", "perfumy dla mężczyzn", "200 ml ", 885, "https://www.notino.pl/paco-rabanne/invictus-victory-elixir-perfumy-dla-mezczyzn/p-16168428/", "11.04.2023"),
    @("Paco Rabanne", "Invictus Victory Elixir ", "perfumy dla mężczyzn", "200 ml ", 885, "https://www.notino.pl/paco-rabanne/invictus-victory-elixir-perfumy-dla-mezczyzn/p-16168428/", "11.04.2023"),
    @("Paco Rabanne", "Invictus Victory Elixir ", "perfumy dla mężczyzn", "200 ml ", 885, "https://www.notino.pl/paco-rabanne/invictus-victory-elixir-perfumy-dla-mezczyzn/p-16168428/", "11.04.2023"),
    @("Parfums", "Roja ", "perfumy unisex", "100 ml ", 10414.8, "https://www.notino.pl/roja-parfums/roja-perfumy-unisex/p-565038/?utm_source=cj&utm_medium=affiliate&utm_campaign=4023395&utm_term=8280252&cjevent=43e4d545d86f11ed821735220a18b8f8&cjdata=MXxZfDB8WXww", "11.04.2023")
)

# Column G holds plain text date strings (e.g. "10.03.2023"), not real
# date values. Force text formatting on the whole new block first so Excel
# doesn't auto-convert the "dd.mm.yyyy" text into a date serial, then strip
# the formatting back off so the cells end up with no explicit style
# (matching the rest of the sheet).
$gRange = $ws.Range("G27:G38")
$gRange.NumberFormat = "@"

$row = 27
foreach ($data in $newRows) {
    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
    $ws.Cells.Item($row, 4).Value = $data[3]
    $ws.Cells.Item($row, 5).Value = $data[4]
    $ws.Cells.Item($row, 6).Value = $data[5]
    $ws.Cells.Item($row, 7).Value = $data[6]
    $row++
}

$gRange.ClearFormats()

$ws.Range("K8").Select() | Out-Null
